$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1 ("treatment")
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Updated values for the existing meta-analysis blocks (rows 3-13)
$ws1.Range("B3").Value = 47.95386330054102
$ws1.Range("C3").Value = 68.68627188929486
$ws1.Range("D3").Value = 95.369720161072

$ws1.Range("B4").Value = 0.2527851075605463
$ws1.Range("C4").Value = 0.6108104510183635
$ws1.Range("D4").Value = 1.124040433708165

$ws1.Range("B5").Value = 0.5112236538172655
$ws1.Range("C5").Value = 0.7946729446679957
$ws1.Range("D5").Value = 1.078018377714276

$ws1.Range("B7").Value = 25.58427801821301
$ws1.Range("C7").Value = 33.00187053526607
$ws1.Range("D7").Value = 41.90138377028548

$ws1.Range("B8").Value = 0.141584780786417
$ws1.Range("C8").Value = 0.3176973392150501
$ws1.Range("D8").Value = 0.5637607160423127

$ws1.Range("B9").Value = 0.3817531520970565
$ws1.Range("C9").Value = 0.5718486946218007
$ws1.Range("D9").Value = 0.761766290314822

$ws1.Range("B11").Value = 6.061409997655358
$ws1.Range("C11").Value = 7.781825747199078
$ws1.Range("D11").Value = 9.831330985074656

$ws1.Range("B12").Value = 0.1215314443590003
$ws1.Range("C12").Value = 0.2870770677167396
$ws1.Range("D12").Value = 0.522547084660171

$ws1.Range("B13").Value = 0.3542269126528542
$ws1.Range("C13").Value = 0.5444230134005977
$ws1.Range("D13").Value = 0.7345137646895203

# New "Speed meta analysis" block (rows 14-17)
$ws1.Range("A14").Value = "Speed meta analysis"

$ws1.Range("A15").Value = "mean (km/day)"
$ws1.Range("B15").Value = 3.755598496599259
$ws1.Range("C15").Value = 4.279469245805568
$ws1.Range("D15").Value = 4.854255501286414

$ws1.Range("A16").Value = "CoV² (RVAR)"
$ws1.Range("B16").Value = 0.03044058476455233
$ws1.Range("C16").Value = 0.0720474140677756
$ws1.Range("D16").Value = 0.1312686343552456

$ws1.Range("A17").Value = "CoV  (RSTD)"
$ws1.Range("B17").Value = 0.1772921414248258
$ws1.Range("C17").Value = 0.2727546661584492
$ws1.Range("D17").Value = 0.3681658259272592

# ---------------------------------------------------------------------------
# Sheet 2 ("control")
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Updated values for the existing meta-analysis blocks (rows 3-13)
$ws2.Range("B3").Value = 47.40142530497118
$ws2.Range("C3").Value = 61.04119260577004
$ws2.Range("D3").Value = 77.3518747829873

$ws2.Range("B4").Value = 0.1042337617494806
$ws2.Range("C4").Value = 0.2573162396158124
$ws2.Range("D4").Value = 0.4783679490412392

$ws2.Range("B5").Value = 0.3284913138669584
$ws2.Range("C5").Value = 0.5161230825894327
$ws2.Range("D5").Value = 0.7037213406888364

$ws2.Range("B7").Value = 20.83790116894635
$ws2.Range("C7").Value = 24.23105480801949
$ws2.Range("D7").Value = 28.01426528797305

$ws2.Range("B8").Value = 0.04076051037919635
$ws2.Range("C8").Value = 0.09547657165538789
$ws2.Range("D8").Value = 0.1730783704576971

$ws2.Range("B9").Value = 0.2050913342341457
$ws2.Range("C9").Value = 0.313888866541604
$ws2.Range("D9").Value = 0.4226188876144621

$ws2.Range("B11").Value = 6.005858961625194
$ws2.Range("C11").Value = 7.50247256957072
$ws2.Range("D11").Value = 9.252154526807171

$ws2.Range("B12").Value = 0.07707815635301282
$ws2.Range("C12").Value = 0.1938706148657898
$ws2.Range("D12").Value = 0.3635900410378016

$ws2.Range("B13").Value = 0.2826414777932059
$ws2.Range("C13").Value = 0.4482561223330648
$ws2.Range("D13").Value = 0.6138697396742772

# New "Speed meta analysis" block (rows 14-17)
$ws2.Range("A14").Value = "Speed meta analysis"

$ws2.Range("A15").Value = "mean (km/day)"
$ws2.Range("B15").Value = 4.07320227638991
$ws2.Range("C15").Value = 4.529512508635405
$ws2.Range("D15").Value = 5.020665245218014

$ws2.Range("A16").Value = "CoV² (RVAR)"
$ws2.Range("B16").Value = 0.01528047553391492
$ws2.Range("C16").Value = 0.04138428929724018
$ws2.Range("D16").Value = 0.08026366661400237

$ws2.Range("A17").Value = "CoV  (RSTD)"
$ws2.Range("B17").Value = 0.126138626415766
$ws2.Range("C17").Value = 0.207585690593424
$ws2.Range("D17").Value = 0.2890940545743764
